# Update Betfair Back/Lay odds for Jogos_do_Dia 2026-01-01 (rows 2-13).
# Cell -> new value map below mirrors the authoritative commit diff;
# every other cell in the sheet is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 1.74
    "K2" = 4.8
    "P2" = 3.1
    "Q2" = 1.45
    "S2" = 2.1
    "T2" = 1.52
    "U2" = 2.8
    "W2" = 2.34
    "X2" = 32
    "Z2" = 46
    "AA2" = 110
    "AE2" = 44
    "AJ2" = 20
    "AN2" = 6
    "AO2" = 28
    "F3" = 1.91
    "G3" = 1.93
    "N3" = 5
    "S3" = 2.62
    "T3" = 1.64
    "U3" = 2.46
    "W3" = 2.06
    "X3" = 22
    "Y3" = 20
    "Z3" = 34
    "AB3" = 12.5
    "AC3" = 9.4
    "AD3" = 17.5
    "AE3" = 44
    "AF3" = 13.5
    "AG3" = 10
    "AH3" = 16
    "AI3" = 44
    "AJ3" = 22
    "AK3" = 18
    "AL3" = 27
    "F4" = 1.77
    "K4" = 4.6
    "Q4" = 1.58
    "N5" = 2.38
    "R5" = 1.15
    "S5" = 3.05
    "T5" = 1.04
    "U5" = 1.04
    "K8" = 980
    "H9" = 3.65
    "I9" = 3.7
    "T9" = 1.89
    "U9" = 2.08
    "W9" = 1.77
    "Z9" = 24
    "AJ9" = 28
    "AN9" = 20
    "H10" = 6.6
    "T10" = 1.76
    "U10" = 2.28
    "W10" = 2.74
    "X10" = 23
    "Y10" = 27
    "AA10" = 170
    "AD10" = 22
    "AG10" = 9.6
    "AK10" = 14.5
    "AO10" = 80
    "H11" = 3.5
    "L11" = 1.21
    "O11" = 1.17
    "T11" = 1.54
    "U11" = 2.5
    "F12" = 2.26
    "G12" = 2.28
    "H12" = 3.35
    "I12" = 3.4
    "J12" = 3.75
    "K12" = 3.8
    "L12" = 1.41
    "N12" = 3.95
    "O12" = 1.31
    "P12" = 2
    "Q12" = 1.98
    "R12" = 1.39
    "S12" = 3.5
    "V12" = 1.41
    "W12" = 1.78
    "Y12" = 13.5
    "Z12" = 23
    "AA12" = 60
    "AB12" = 10.5
    "AC12" = 8
    "AD12" = 14.5
    "AE12" = 40
    "AF12" = 14
    "AI12" = 48
    "AL12" = 36
    "AN12" = 17.5
    "AO12" = 36
    "R13" = 1.54
    "AB13" = 32
    "AH13" = 24
    "AJ13" = 290
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
